# US1377 - Fixed some actions from ComplianceReportActions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Users" (sheet1.xml) - move selection to B5
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B5").Select()

# ---------------------------------------------------------------------------
# Sheet "Driver View Test Data" (sheet3.xml) - move selection to D19
# ---------------------------------------------------------------------------
$wsDriver = $wb.Worksheets.Item("Driver View Test Data")
$wsDriver.Range("D19").Select()

# ---------------------------------------------------------------------------
# Sheet "Compliance Report Test Data" (sheet4.xml)
# ---------------------------------------------------------------------------
$wsCompliance = $wb.Worksheets.Item("Compliance Report Test Data")

# -- Edit existing row 2 --
$wsCompliance.Range("C2").Value = 1
$wsCompliance.Range("F2").ClearContents()
$wsCompliance.Range("S2").Value = $false
$wsCompliance.Range("U2").Value = " "
$wsCompliance.Range("V2").Value = " "
$wsCompliance.Range("Y2").Value = "5"
$wsCompliance.Range("Z2").Value = 1

# -- New row 3 --
$wsCompliance.Range("A3").Value = 2
$wsCompliance.Range("B3").Value = "GenerateRandomString(10)"
$wsCompliance.Range("C3").Value = 1
$wsCompliance.Range("D3").Value = "Eastern Time (US and Canada)"
$wsCompliance.Range("E3").Value = 0
$wsCompliance.Range("K3").Value = "District"
$wsCompliance.Range("L3").Value = "<SPECIFY_BOUNDARY_NAME>"
$wsCompliance.Range("P2:Q2").Copy()
$wsCompliance.Range("P3:Q3").PasteSpecial(-4122)
$wsCompliance.Range("P3").Value = 42344.500011574077
$wsCompliance.Range("Q3").Value = 42010.500011574077
$wsCompliance.Range("R3").Value = "Rapid Response"
$wsCompliance.Range("W3").Value = 8.5
$wsCompliance.Range("X3").Value = 11
$wsCompliance.Range("Y3").Value = "6,7,8"
$wsCompliance.Range("Z3").Value = 5
$wsCompliance.Range("AA3").Value = 2

# -- New row 4 --
$wsCompliance.Range("A4").Value = 3
$wsCompliance.Range("B4").Value = "GenerateRandomString(10)"
$wsCompliance.Range("C4").Value = 1
$wsCompliance.Range("D4").Value = "Mountain Time (US and Canada)"
$wsCompliance.Range("G4").Value = 37.435339792682498
$wsCompliance.Range("H4").Value = -121.846961975097
$wsCompliance.Range("I4").Value = 37.330583620739603
$wsCompliance.Range("J4").Value = -122.04883575439401
$wsCompliance.Range("M2:N2").Copy()
$wsCompliance.Range("M4:N4").PasteSpecial(-4122)
$wsCompliance.Range("O4").Value = "stnd-sqacudr"
$wsCompliance.Range("P2:Q2").Copy()
$wsCompliance.Range("P4:Q4").PasteSpecial(-4122)
$wsCompliance.Range("P4").Value = 42294.500011574077
$wsCompliance.Range("Q4").Value = 42010.500011574077
$wsCompliance.Range("R4").Value = "Standard"
$wsCompliance.Range("S4").Value = $false
$wsCompliance.Range("T4").Value = 2
$wsCompliance.Range("U4").Value = " "
$wsCompliance.Range("V4").Value = " "
$wsCompliance.Range("W4").Value = 8.5
$wsCompliance.Range("X4").Value = 11
$wsCompliance.Range("Y4").Value = "5"
$wsCompliance.Range("Z4").Value = 1
$wsCompliance.Range("AA4").Value = 4

# Y2/Y4 keep the "text" number format already present on Y2; reapply to Y4
$wsCompliance.Range("Y2").Copy()
$wsCompliance.Range("Y4").PasteSpecial(-4122)
$wsCompliance.Range("Y4").Value = "5"

# Restore the selection / view for this sheet
$wsCompliance.Range("Y4").Select()

# ---------------------------------------------------------------------------
# Sheet "Report Views Data" (sheet5.xml)
# ---------------------------------------------------------------------------
$wsViews = $wb.Worksheets.Item("Report Views Data")
$wsViews.Range("N5").Value = "All options selected (except IsoCap and Annotation), BaseMap=NONE"

$wsViews.Range("A6").Value = 5
$wsViews.Range("B6").Value = "GenerateRandomString(10)"
$wsViews.Range("C6").Value = $false
$wsViews.Range("D6").Value = $false
$wsViews.Range("E6").Value = $false
$wsViews.Range("F6").Value = $true
$wsViews.Range("G6").Value = $true
$wsViews.Range("H6").Value = $false
$wsViews.Range("I6").Value = $true
$wsViews.Range("J6").Value = $true
$wsViews.Range("K6").Value = $false
$wsViews.Range("L6").Value = $false
$wsViews.Range("M6").Value = "Indication Table"

$wsViews.Range("A7").Value = 6
$wsViews.Range("B7").Value = "GenerateRandomString(10)"
$wsViews.Range("C7").Value = $true
$wsViews.Range("D7").Value = $true
$wsViews.Range("E7").Value = $false
$wsViews.Range("F7").Value = $true
$wsViews.Range("G7").Value = $false
$wsViews.Range("H7").Value = $false
$wsViews.Range("I7").Value = $true
$wsViews.Range("J7").Value = $true
$wsViews.Range("K7").Value = $false
$wsViews.Range("L7").Value = $false
$wsViews.Range("M7").Value = "Gap Table"

$wsViews.Range("A8").Value = 7
$wsViews.Range("B8").Value = "GenerateRandomString(10)"
$wsViews.Range("C8").Value = $false
$wsViews.Range("D8").Value = $false
$wsViews.Range("E8").Value = $true
$wsViews.Range("F8").Value = $true
$wsViews.Range("G8").Value = $false
$wsViews.Range("H8").Value = $false
$wsViews.Range("I8").Value = $true
$wsViews.Range("J8").Value = $true
$wsViews.Range("K8").Value = $false
$wsViews.Range("L8").Value = $false
$wsViews.Range("M8").Value = "Indication Table"

$wsViews.Range("A9").Value = 8
$wsViews.Range("B9").Value = "GenerateRandomString(10)"
$wsViews.Range("C9").Value = $false
$wsViews.Range("D9").Value = $true
$wsViews.Range("E9").Value = $false
$wsViews.Range("F9").Value = $false
$wsViews.Range("G9").Value = $false
$wsViews.Range("H9").Value = $true
$wsViews.Range("I9").Value = $false
$wsViews.Range("J9").Value = $false
$wsViews.Range("K9").Value = $false
$wsViews.Range("L9").Value = $false
$wsViews.Range("M9").Value = "Gap Table"

$wsViews.Range("A10").Select()

# ---------------------------------------------------------------------------
# Sheet "Report Opt View Layers" (sheet6.xml)
# ---------------------------------------------------------------------------
$wsLayers = $wb.Worksheets.Item("Report Opt View Layers")

$wsLayers.Range("A6").Value = 5
$wsLayers.Range("B6").Value = $false
$wsLayers.Range("C6").Value = $true
$wsLayers.Range("D6").Value = $true
$wsLayers.Range("E6").Value = $true
$wsLayers.Range("F6").Value = $false
$wsLayers.Range("G6").Value = $false
$wsLayers.Range("H6").Value = $false
$wsLayers.Range("I6").Value = $false

# Remove the now-unused trailing empty formatted row (mirrors the deleted
# blank row 602 in the original workbook, which shifted everything up by one
# and dropped the final blank row 736).
$wsLayers.Rows("602:602").Delete()

$wsLayers.Range("E7").Select()

# ---------------------------------------------------------------------------
# Sheet "Report Opt Tabular PDF Content" (sheet7.xml)
# ---------------------------------------------------------------------------
$wsTabular = $wb.Worksheets.Item("Report Opt Tabular PDF Content")
$wsTabular.Rows("579:579").Delete()
$wsTabular.Rows("7:7").Select()

# ---------------------------------------------------------------------------
# Re-activate the workbook's primary sheet last, so "tabSelected" stays on
# "Compliance Report Test Data" like in the source file.
# ---------------------------------------------------------------------------
$wsCompliance.Activate()
$wsCompliance.Range("Y4").Select()
